# The deck ships two theme parts:
#   theme1.xml ("Office Theme" colours)  -- linked from the Notes Master
#   theme2.xml ("Integral" colours)      -- linked from the Slide Master
#   (the Slide Master's theme is what actually paints the slides)
#
# The target revision swaps which colour palette each master uses: the
# Slide Master's theme should end up holding the "Office Theme" palette
# (it previously held "Integral"). We apply that by pushing the
# "Office Theme" colour values onto the presentation's live colour
# scheme (PowerPoint resolves ColorScheme/ThemeColorScheme against the
# Slide Master's theme part, i.e. theme2.xml).

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

function RGBInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2),16)
    $g = [Convert]::ToInt32($hex.Substring(2,2),16)
    $b = [Convert]::ToInt32($hex.Substring(4,2),16)
    return $r + ($g * 256) + ($b * 65536)
}

# "Office Theme" palette: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeTheme = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$cs = $s1.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $cs.Colors($i).RGB = RGBInt($officeTheme[$i - 1])
}

# Best-effort: keep the theme / colour-scheme display names in sync with
# the newly-applied palette (no-op on hosts that don't support renaming
# the theme part, but harmless).
try {
    $p.Designs.Item(1).Name = "Office Theme"
} catch {}
try {
    $s1.ThemeColorScheme.Name = "Office"
} catch {}
try {
    $p.SlideMaster.Theme.Name = "Office Theme"
} catch {}
